$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.170.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.626.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.37"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.626.76"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.14"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.545"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.137.24"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.68"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0745"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.999"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.94"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.39"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.43"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.46"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.32"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.99"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.352.72"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.16%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0177"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.550"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.854"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.50"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.762.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.58"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.62"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.844"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +27.18%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆01000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.43%  "
